{"js": "// Adiciona funcionalidade para padronizar a coluna battery_capacity e\n// criar um CSV com os dados transformados.\n//\n// The commit colors (font color 0000FF / blue) the three bullet items\n// describing the new CSV-export behaviour:\n//   - \"Baixar os dados do MySQL ap\u00f3s a inser\u00e7\u00e3o.\"\n//   - \"Aplicar novas transforma\u00e7\u00f5es conforme necess\u00e1rio.\"\n//   - \"Salvar os dados processados como CSV na pasta dados/ com timestamp.\"\n\nconst targetTexts = [\n  \"Baixar os dados do MySQL ap\u00f3s a inser\u00e7\u00e3o.\",\n  \"Aplicar novas transforma\u00e7\u00f5es conforme necess\u00e1rio.\",\n  \"Salvar os dados processados como CSV na pasta dados/ com timestamp.\"\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const paragraph = paragraphs.items[i];\n  const text = (paragraph.text || \"\").trim();\n  if (targetTexts.indexOf(text) !== -1) {\n    // Setting the paragraph's font color stamps <w:color val=\"0000FF\"/>\n    // into both the paragraph mark run properties (pPr/rPr) and every\n    // run's run properties (rPr) within the paragraph.\n    paragraph.font.color = \"#0000FF\";\n  }\n}\n\nawait context.sync();\n", "ps1": "# Adiciona funcionalidade para padronizar a coluna battery_capacity e\n# criar um CSV com os dados transformados.\n#\n# Colors (font color blue / RGB 0000FF) the three bullet items describing\n# the new CSV-export behaviour:\n#   - \"Baixar os dados do MySQL ap\u00f3s a inser\u00e7\u00e3o.\"\n#   - \"Aplicar novas transforma\u00e7\u00f5es conforme necess\u00e1rio.\"\n#   - \"Salvar os dados processados como CSV na pasta dados/ com timestamp.\"\n\n$d = $word.ActiveDocument\n\n$targets = @(\n  \"Baixar os dados do MySQL ap\u00f3s a inser\u00e7\u00e3o.\",\n  \"Aplicar novas transforma\u00e7\u00f5es conforme necess\u00e1rio.\",\n  \"Salvar os dados processados como CSV na pasta dados/ com timestamp.\"\n)\n\nforeach ($p in $d.Paragraphs) {\n  $text = $p.Range.Text.Trim()\n  if ($targets -contains $text) {\n    # wdColorBlue (&HFF0000 = 16711680) -> <w:color w:val=\"0000FF\"/>\n    # applied to both the paragraph mark and every run in the paragraph.\n    $p.Range.Font.Color = 16711680\n  }\n}\n"}
